$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.025.25'
$ws.Range("E2").Value = '  -0.88%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.827.37'
$ws.Range("E3").Value = '  -0.88%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("E5").Value = '  +0.22%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6362'
$ws.Range("E6").Value = '  -5.33%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.59'
$ws.Range("E8").Value = '  +5.82%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07349'
$ws.Range("E9").Value = '  -1.20%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2932'
$ws.Range("E10").Value = '  -0.25%  '

$ws.Range("E11").Value = '  +0.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07669'
$ws.Range("E12").Value = '  -0.67%  '

$ws.Range("E13").Value = '  -0.67%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.994'
$ws.Range("E14").Value = '  -0.31%  '

$ws.Range("E15").Value = '  -0.93%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.24'
$ws.Range("E16").Value = '  -4.41%  '

$ws.Range("E17").Value = '  -1.45%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008629'
$ws.Range("E18").Value = '  +3.70%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.031.22'
$ws.Range("E19").Value = '  -0.64%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.081.06'
$ws.Range("E20").Value = '  +0.83%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.42'
$ws.Range("E21").Value = '  -0.82%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '224.12'
$ws.Range("E22").Value = '  -2.00%  '

$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.109'
$ws.Range("E24").Value = '  -0.74%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.001'
$ws.Range("E25").Value = '  +0.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.18'
$ws.Range("E26").Value = '  -1.71%  '

$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1380'
$ws.Range("E27").Value = '  -1.76%  '

$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.462'
$ws.Range("E28").Value = '  -2.80%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '17.90'
$ws.Range("E29").Value = '  -0.72%  '

$ws.Range("E30").Value = '  -0.79%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.090'
$ws.Range("E31").Value = '  -1.73%  '

$ws.Range("E33").Value = '  +0.51%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05300'
$ws.Range("E34").Value = '  +0.19%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7421'
$ws.Range("E35").Value = '  -1.16%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.831'
$ws.Range("E36").Value = '  -2.27%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.152'
$ws.Range("E37").Value = '  +1.35%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.647'
$ws.Range("E38").Value = '  -1.20%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.288.70'
$ws.Range("E39").Value = '  -2.44%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01782'
$ws.Range("E40").Value = '  -1.31%  '

$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.738'
$ws.Range("E41").Value = '  +0.43%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.350'
$ws.Range("E42").Value = '  +6.35%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8970'
$ws.Range("E43").Value = '  -2.51%  '

$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9997'
$ws.Range("E44").Value = '  -0.24%  '

$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.79'
$ws.Range("E45").Value = '  +0.63%  '

$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000125'
$ws.Range("E46").Value = '  +3.35%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.978.12'
$ws.Range("E47").Value = '  +0.41%  '

$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '64.23'
$ws.Range("E48").Value = '  +0.61%  '

$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.5136'
$ws.Range("E49").Value = '  -0.55%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.731'
$ws.Range("E50").Value = '  -2.49%  '

$ws.Range("E51").Value = '  -1.95%  '
